# Refresh the NATMI Angpt1-Tie1 LR-pair table with the new TPM-based run.
# Rows 2-13 get updated metrics (a 4th cluster, "Resolving-Mac", was added to
# the analysis), and rows 14-17 are newly appended for "Resolving-Mac" as a
# sending cluster.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{}
$data[2] = @("ECs", "Angpt1", "Tie1", "ECs", 2, 0.6666666666666666, 0.1786683333333333, 0.536005, 0.005700931134181372, 0.005700931134181372, 3, 1, 75.02398799999999, 225.071964, 0.9946207163417996, 0.9946207163417997, 13.40441089598, 120.63969806382, 0.005670264208494744, 0.005670264208494745)
$data[3] = @("ECs", "Angpt1", "Tie1", "FAPs", 2, 0.6666666666666666, 0.1786683333333333, 0.536005, 0.005700931134181372, 0.005700931134181372, 1, 0.3333333333333333, 0.09667466666666667, 0.290024, 0.001281651759329359, 0.001281651759329359, 0.01727270156888889, 0.15545431412, 0.000007306608417939073, 0.000007306608417939073)
$data[4] = @("ECs", "Angpt1", "Tie1", "MuSCs", 2, 0.6666666666666666, 0.1786683333333333, 0.536005, 0.005700931134181372, 0.005700931134181372, 3, 1, 0.291865, 0.875595, 0.00386936209489556, 0.00386936209489556, 0.05214703310833332, 0.4693232979749999, 0.00002205896683621135, 0.00002205896683621136)
$data[5] = @("ECs", "Angpt1", "Tie1", "Resolving-Mac", 2, 0.6666666666666666, 0.1786683333333333, 0.536005, 0.005700931134181372, 0.005700931134181372, 1, 0.3333333333333333, 0.01721833333333333, 0.051655, 0.0002282698039753883, 0.0002282698039753884, 0.003076370919444444, 0.027687338275, 0.00000130135043247677, 0.00000130135043247677)
$data[6] = @("FAPs", "Angpt1", "Tie1", "ECs", 3, 1, 22.14783133333333, 66.443494, 0.7066907652137446, 0.7066907652137446, 3, 1, 75.02398799999999, 225.071964, 0.9946207163417996, 0.9946207163417997, 1661.618632178024, 14954.56768960221, 0.7028892751290292, 0.7028892751290293)
$data[7] = @("FAPs", "Angpt1", "Tie1", "FAPs", 3, 1, 22.14783133333333, 66.443494, 0.7066907652137446, 0.7066907652137446, 1, 0.3333333333333333, 0.09667466666666667, 0.290024, 0.001281651759329359, 0.001281651759329359, 2.141134211539556, 19.270207903856, 0.0009057314625380067, 0.0009057314625380067)
$data[8] = @("FAPs", "Angpt1", "Tie1", "MuSCs", 3, 1, 22.14783133333333, 66.443494, 0.7066907652137446, 0.7066907652137446, 3, 1, 0.291865, 0.875595, 0.00386936209489556, 0.00386936209489556, 6.464176792103332, 58.17759112893, 0.002734442459730801, 0.002734442459730802)
$data[9] = @("FAPs", "Angpt1", "Tie1", "Resolving-Mac", 3, 1, 22.14783133333333, 66.443494, 0.7066907652137446, 0.7066907652137446, 1, 0.3333333333333333, 0.01721833333333333, 0.051655, 0.0002282698039753883, 0.0002282698039753884, 0.3813487425077777, 3.43213868257, 0.0001613161624465587, 0.0001613161624465587)
$data[10] = @("MuSCs", "Angpt1", "Tie1", "ECs", 3, 1, 8.967326666666667, 26.90198, 0.2861285535640985, 0.2861285535640984, 3, 1, 75.02398799999999, 225.071964, 0.9946207163417996, 0.9946207163417997, 672.7646082320799, 6054.88147408872, 0.2845893869117666, 0.2845893869117666)
$data[11] = @("MuSCs", "Angpt1", "Tie1", "FAPs", 3, 1, 8.967326666666667, 26.90198, 0.2861285535640985, 0.2861285535640984, 1, 0.3333333333333333, 0.09667466666666667, 0.290024, 0.001281651759329359, 0.001281651759329359, 0.8669133163911111, 7.802219847520001, 0.0003667171640697916, 0.0003667171640697915)
$data[12] = @("MuSCs", "Angpt1", "Tie1", "MuSCs", 3, 1, 8.967326666666667, 26.90198, 0.2861285535640985, 0.2861285535640984, 3, 1, 0.291865, 0.875595, 0.00386936209489556, 0.00386936209489556, 2.617248797566667, 23.5552391781, 0.001107134979428217, 0.001107134979428217)
$data[13] = @("MuSCs", "Angpt1", "Tie1", "Resolving-Mac", 3, 1, 8.967326666666667, 26.90198, 0.2861285535640985, 0.2861285535640984, 1, 0.3333333333333333, 0.01721833333333333, 0.051655, 0.0002282698039753883, 0.0002282698039753884, 0.1544024196555555, 1.3896217769, 0.00006531450883383817, 0.00006531450883383816)
$data[14] = @("Resolving-Mac", "Angpt1", "Tie1", "ECs", 2, 0.6666666666666666, 0.04637566666666667, 0.139127, 0.001479750087975396, 0.001479750087975395, 3, 1, 75.02398799999999, 225.071964, 0.9946207163417996, 0.9946207163417997, 3.479287459492, 31.313587135428, 0.001471790092508929, 0.001471790092508929)
$data[15] = @("Resolving-Mac", "Angpt1", "Tie1", "FAPs", 2, 0.6666666666666666, 0.04637566666666667, 0.139127, 0.001479750087975396, 0.001479750087975395, 1, 0.3333333333333333, 0.09667466666666667, 0.290024, 0.001281651759329359, 0.001281651759329359, 0.004483352116444445, 0.040350169048, 0.000001896524303621439, 0.000001896524303621439)
$data[16] = @("Resolving-Mac", "Angpt1", "Tie1", "MuSCs", 2, 0.6666666666666666, 0.04637566666666667, 0.139127, 0.001479750087975396, 0.001479750087975395, 3, 1, 0.291865, 0.875595, 0.00386936209489556, 0.00386936209489556, 0.01353543395166667, 0.121818905565, 0.000005725688900330366, 0.000005725688900330366)
$data[17] = @("Resolving-Mac", "Angpt1", "Tie1", "Resolving-Mac", 2, 0.6666666666666666, 0.04637566666666667, 0.139127, 0.001479750087975396, 0.001479750087975395, 1, 0.3333333333333333, 0.01721833333333333, 0.051655, 0.0002282698039753883, 0.0002282698039753884, 0.0007985116872222222, 0.007186605185, 0.0000003377822625147072, 0.0000003377822625147072)

foreach ($r in $data.Keys) {
    $row = $data[$r]
    for ($c = 1; $c -le $row.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}

Write-Output "done"